$d = $word.ActiveDocument

# --- Paragraph describing how the data was extracted from NBA.com ---
# (was: "We extracted the data from NBA.com APIs ..."
#  now: "I extracted my data from NBA.com ...")
$range = $d.Paragraphs.Item(9).Range
$range.Find.Execute(
    "We extracted the data from NBA.com APIs with the help from",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I extracted my data from NBA.com with the help from", 2)

$range = $d.Paragraphs.Item(9).Range
$range.Find.Execute(
    "With their extensive documentation and example on how to extract the data from the various endpoints",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "They provided extensive documentation and example on to use their client package to extract from the various endpoints",
    2)

$range = $d.Paragraphs.Item(9).Range
$range.Find.Execute(
    "I created a Postgres SQL Database to store our extracted data into 5 different tables.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Once the data was successfully extracted, I proceeded to load the data into Postgres SQL Database to store.",
    2)

$range = $d.Paragraphs.Item(9).Range
$range.Find.Execute(
    "As you can see from my ERD we have a teams,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As you can see from my ERD  we have 5 tables, teams,", 2)

# --- Model-performance paragraph: tidy up "in order to ..." phrase ---
$range = $d.Paragraphs.Item(16).Range
$range.Find.Execute(
    "in order to determine how well our models really performed, what",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in order to determine how well our models really performed, what",
    2)

Write-Output $d.Paragraphs.Item(9).Range.Text
Write-Output $d.Paragraphs.Item(16).Range.Text
